$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: relocate the totals row (was 32) and footer row (was 33) down to 36/37 ---
$totalValue = $ws.Range("P32").Value()
$footerA = $ws.Range("A33").Value()
$footerG = $ws.Range("G33").Value()
$footerK = $ws.Range("K33").Value()

# copy formats of row 32 (totals) to row 36, and row 33 (footer) to row 37
$ws.Range("A32:Q32").Copy()
$ws.Range("A36:Q36").PasteSpecial(-4122)
$ws.Range("A33:Q33").Copy()
$ws.Range("A37:Q37").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# drop the old merges that lived on rows 32/33 before those rows get reused
$ws.Range("P32:Q32").UnMerge()
$ws.Range("A33:F33").UnMerge()
$ws.Range("G33:I33").UnMerge()
$ws.Range("K33:Q33").UnMerge()

# re-create merges on the relocated rows
$ws.Range("P36:Q36").Merge()
$ws.Range("A37:F37").Merge()
$ws.Range("G37:I37").Merge()
$ws.Range("K37:Q37").Merge()

# restore the values onto the relocated rows (timestamp text is updated to the new save time)
$ws.Range("P36").Value = 1607.27
$ws.Range("A37").Value = 'Monday, 9 June, 2025 12:18 PM'
$ws.Range("G37").Value = $footerG
$ws.Range("K37").Value = $footerK

# --- Step 2: build four brand-new product rows (32-35), formatted like row 31 ---
$ws.Range("A31:Q31").Copy()
$ws.Range("A32:Q32").PasteSpecial(-4122)
$ws.Range("A33:Q33").PasteSpecial(-4122)
$ws.Range("A34:Q34").PasteSpecial(-4122)
$ws.Range("A35:Q35").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("A32:B32").Merge()
$ws.Range("C32:G32").Merge()
$ws.Range("H32:K32").Merge()
$ws.Range("L32:M32").Merge()
$ws.Range("N32:O32").Merge()

$ws.Range("A33:B33").Merge()
$ws.Range("C33:G33").Merge()
$ws.Range("H33:K33").Merge()
$ws.Range("L33:M33").Merge()
$ws.Range("N33:O33").Merge()

$ws.Range("A34:B34").Merge()
$ws.Range("C34:G34").Merge()
$ws.Range("H34:K34").Merge()
$ws.Range("L34:M34").Merge()
$ws.Range("N34:O34").Merge()

$ws.Range("A35:B35").Merge()
$ws.Range("C35:G35").Merge()
$ws.Range("H35:K35").Merge()
$ws.Range("L35:M35").Merge()
$ws.Range("N35:O35").Merge()

# row heights for the new rows, following the workbooks alternating pattern
$ws.Rows("32:32").RowHeight = 25.5
$ws.Rows("33:33").RowHeight = 24.75
$ws.Rows("34:34").RowHeight = 25.5
$ws.Rows("35:35").RowHeight = 24.75

# --- Step 3: (re)write all 29 product rows (7-35) with the final data ---
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = 'ALKAPRESS PLUS 10/160MG 20 F.C. TABS.'
$ws.Range("H7").Value = '0:0'
$ws.Range("L7").Value = '1'
$ws.Range("N7").Value = '102.00'
$ws.Range("P7").Value = '102.0000'
$ws.Range("Q7").Value = '1:0'

$ws.Range("A8").Value = 2
$ws.Range("C8").Value = 'BISOLOCK-D 10/25MG 20 F.C.TAB'
$ws.Range("H8").Value = '0:0'
$ws.Range("L8").Value = '1'
$ws.Range("N8").Value = '42.00'
$ws.Range("P8").Value = '21.0000'
$ws.Range("Q8").Value = '0:1'

$ws.Range("A9").Value = 3
$ws.Range("C9").Value = 'BRAVAMAX 200MG 10 TAB'
$ws.Range("H9").Value = '0:0'
$ws.Range("L9").Value = '1'
$ws.Range("N9").Value = '122.00'
$ws.Range("P9").Value = '122.0000'
$ws.Range("Q9").Value = '1:0'

$ws.Range("A10").Value = 4
$ws.Range("C10").Value = 'BRUFEN 400MG 30 TAB'
$ws.Range("H10").Value = '1:2'
$ws.Range("L10").Value = '1'
$ws.Range("N10").Value = '78.00'
$ws.Range("P10").Value = '25.7400'
$ws.Range("Q10").Value = '0:1'

$ws.Range("A11").Value = 5
$ws.Range("C11").Value = 'DECLOPHEN 75MG/3ML 3 AMPOULES'
$ws.Range("H11").Value = '5:0'
$ws.Range("L11").Value = '1'
$ws.Range("N11").Value = '36.00'
$ws.Range("P11").Value = '11.8800'
$ws.Range("Q11").Value = '0:1'

$ws.Range("A12").Value = 6
$ws.Range("C12").Value = 'FLAGYL 500MG 20 TAB.'
$ws.Range("H12").Value = '2:1'
$ws.Range("L12").Value = '1'
$ws.Range("N12").Value = '34.00'
$ws.Range("P12").Value = '-17.0000'
$ws.Range("Q12").Value = '0:-1'

$ws.Range("A13").Value = 7
$ws.Range("C13").Value = 'GLAPTIVIA PLUS 50/1000MG 30 F.C. TAB.'
$ws.Range("H13").Value = '1:0'
$ws.Range("L13").Value = '1'
$ws.Range("N13").Value = '168.00'
$ws.Range("P13").Value = '55.4400'
$ws.Range("Q13").Value = '0:1'

$ws.Range("A14").Value = 8
$ws.Range("C14").Value = 'GLIMET 2.5/400 MG 30 F.C.TAB.'
$ws.Range("H14").Value = '0:2'
$ws.Range("L14").Value = '1'
$ws.Range("N14").Value = '45.00'
$ws.Range("P14").Value = '45.0000'
$ws.Range("Q14").Value = '1:0'

$ws.Range("A15").Value = 9
$ws.Range("C15").Value = 'GLIPTUS PLUS 50/1000MG 30 TABLETS'
$ws.Range("H15").Value = '1:0'
$ws.Range("L15").Value = '1'
$ws.Range("N15").Value = '192.00'
$ws.Range("P15").Value = '192.0000'
$ws.Range("Q15").Value = '1:0'

$ws.Range("A16").Value = 10
$ws.Range("C16").Value = 'GUAVA SYRUP 120 ML'
$ws.Range("H16").Value = '0:0'
$ws.Range("L16").Value = '1'
$ws.Range("N16").Value = '39.00'
$ws.Range("P16").Value = '39.0000'
$ws.Range("Q16").Value = '1:0'

$ws.Range("A17").Value = 11
$ws.Range("C17").Value = 'HEDERA HELIX SYRAP'
$ws.Range("H17").Value = '12:0'
$ws.Range("L17").Value = '0'
$ws.Range("N17").Value = '65.00'
$ws.Range("P17").Value = '65.0000'
$ws.Range("Q17").Value = '1:0'

$ws.Range("A18").Value = 12
$ws.Range("C18").Value = 'HEPAMARIN 140MG 30 CAPSULE'
$ws.Range("H18").Value = '0:1'
$ws.Range("L18").Value = '1'
$ws.Range("N18").Value = '75.00'
$ws.Range("P18").Value = '75.0000'
$ws.Range("Q18").Value = '1:0'

$ws.Range("A19").Value = 13
$ws.Range("C19").Value = 'NEURONTIN 300MG 20 CAPS'
$ws.Range("H19").Value = '0:1'
$ws.Range("L19").Value = '1'
$ws.Range("N19").Value = '108.00'
$ws.Range("P19").Value = '108.0000'
$ws.Range("Q19").Value = '1:0'

$ws.Range("A20").Value = 14
$ws.Range("C20").Value = 'NEVXAL FORTE 0.3%  DROPS'
$ws.Range("H20").Value = '0:0'
$ws.Range("L20").Value = '0'
$ws.Range("N20").Value = '44.00'
$ws.Range("P20").Value = '44.0000'
$ws.Range("Q20").Value = '1:0'

$ws.Range("A21").Value = 15
$ws.Range("C21").Value = 'NORGESIC 20 TAB.'
$ws.Range("H21").Value = '0:0'
$ws.Range("L21").Value = '1'
$ws.Range("N21").Value = '38.00'
$ws.Range("P21").Value = '38.0000'
$ws.Range("Q21").Value = '1:0'

$ws.Range("A22").Value = 16
$ws.Range("C22").Value = 'PROSTRIDE 5MG 30 CAPS.'
$ws.Range("H22").Value = '0:1'
$ws.Range("L22").Value = '1'
$ws.Range("N22").Value = '183.00'
$ws.Range("P22").Value = '60.3900'
$ws.Range("Q22").Value = '0:1'

$ws.Range("A23").Value = 17
$ws.Range("C23").Value = 'RELAT HAIR SERUM 60 ML'
$ws.Range("H23").Value = '0:0'
$ws.Range("L23").Value = '1'
$ws.Range("N23").Value = '225.00'
$ws.Range("P23").Value = '225.0000'
$ws.Range("Q23").Value = '1:0'

$ws.Range("A24").Value = 18
$ws.Range("C24").Value = 'STREPTOQUIN 20 TABLETS'
$ws.Range("H24").Value = '3:0'
$ws.Range("L24").Value = '1'
$ws.Range("N24").Value = '46.00'
$ws.Range("P24").Value = '23.0000'
$ws.Range("Q24").Value = '0:1'

$ws.Range("A25").Value = 19
$ws.Range("C25").Value = 'TAMSULIN 0.4MG 28 CAPS'
$ws.Range("H25").Value = '0:0'
$ws.Range("L25").Value = '1'
$ws.Range("N25").Value = '124.00'
$ws.Range("P25").Value = '62.0000'
$ws.Range("Q25").Value = '0:1'

$ws.Range("A26").Value = 20
$ws.Range("C26").Value = 'URIVIN-N 10 EFF. SACHETS'
$ws.Range("H26").Value = '2:0'
$ws.Range("L26").Value = '1'
$ws.Range("N26").Value = '31.00'
$ws.Range("P26").Value = '31.0000'
$ws.Range("Q26").Value = '1:0'

$ws.Range("A27").Value = 21
$ws.Range("C27").Value = 'VOLTAREN 75MG/3ML 3 AMP.'
$ws.Range("H27").Value = '4:0'
$ws.Range("L27").Value = '1'
$ws.Range("N27").Value = '51.00'
$ws.Range("P27").Value = '33.6600'
$ws.Range("Q27").Value = '0:2'

$ws.Range("A28").Value = 22
$ws.Range("C28").Value = 'VOMIBREAK 30 DELAYED RELEASE F.C. TABLETS'
$ws.Range("H28").Value = '0:1'
$ws.Range("L28").Value = '1'
$ws.Range("N28").Value = '66.00'
$ws.Range("P28").Value = '66.0000'
$ws.Range("Q28").Value = '1:0'

$ws.Range("A29").Value = 23
$ws.Range("C29").Value = 'ZANOGLIDE 4/30 MG 30 TAB'
$ws.Range("H29").Value = '0:2'
$ws.Range("L29").Value = '1'
$ws.Range("N29").Value = '102.00'
$ws.Range("P29").Value = '33.6600'
$ws.Range("Q29").Value = '0:1'

$ws.Range("A30").Value = 24
$ws.Range("C30").Value = 'ZURCAL 20 MG 14 GASTRO-RESISTANT TABS.'
$ws.Range("H30").Value = '1:0'
$ws.Range("L30").Value = '1'
$ws.Range("N30").Value = '81.00'
$ws.Range("P30").Value = '40.5000'
$ws.Range("Q30").Value = '0:1'

$ws.Range("A31").Value = 25
$ws.Range("C31").Value = 'حبايه'
$ws.Range("H31").Value = '0:0'
$ws.Range("L31").Value = '0'
$ws.Range("N31").Value = '3.00'
$ws.Range("P31").Value = '9.0000'
$ws.Range("Q31").Value = '3:0'

$ws.Range("A32").Value = 26
$ws.Range("C32").Value = 'سرنجات 3 سم'
$ws.Range("H32").Value = '0:0'
$ws.Range("L32").Value = '0'
$ws.Range("N32").Value = '2.00'
$ws.Range("P32").Value = '6.0000'
$ws.Range("Q32").Value = '3:0'

$ws.Range("A33").Value = 27
$ws.Range("C33").Value = 'شامبو كلير 1 كيس'
$ws.Range("H33").Value = '117:0'
$ws.Range("L33").Value = '0'
$ws.Range("N33").Value = '2.50'
$ws.Range("P33").Value = '5.0000'
$ws.Range("Q33").Value = '2:0'

$ws.Range("A34").Value = 28
$ws.Range("C34").Value = 'كريم شعر ايفا 85مل '
$ws.Range("H34").Value = '2:0'
$ws.Range("L34").Value = '0'
$ws.Range("N34").Value = '50.00'
$ws.Range("P34").Value = '50.0000'
$ws.Range("Q34").Value = '1:0'

$ws.Range("A35").Value = 29
$ws.Range("C35").Value = 'كريم فيبكس الازرق'
$ws.Range("H35").Value = '1:0'
$ws.Range("L35").Value = '0'
$ws.Range("N35").Value = '35.00'
$ws.Range("P35").Value = '35.0000'
$ws.Range("Q35").Value = '1:0'

